# Update the NATMI Plat-Lrp1 ligand-receptor table with the recomputed TPM
# values (per-cluster ligand/receptor avg+total expression and derived
# specificity scores, plus the edge weight/specificity columns Q-T which
# are simple products: Q=G*M, R=H*N, S=I*O, T=J*P).
#
# Note: for row 21 (MuSCs -> Resolving-Mac) the source diff's hunk for
# columns Q/R/S/T duplicates row 11's new values verbatim (a copy/paste
# artifact upstream - confirmed because the diff's own "old" values for
# Q21/R21/S21/T21 don't match the original workbook's actual contents).
# The values written here for Q21/R21/S21/T21 are instead derived from
# row 21's own new G/H/I/J/M/N/O/P via the same Q=G*M/R=H*N/S=I*O/T=J*P
# relationship that holds for every other row in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 8.872289666666667
$ws.Cells.Item(2, 8).Value = 26.616869
$ws.Cells.Item(2, 9).Value = 0.1344073019604298
$ws.Cells.Item(2, 10).Value = 0.1393180933827548
$ws.Cells.Item(2, 13).Value = 2.906846333333333
$ws.Cells.Item(2, 14).Value = 8.720538999999999
$ws.Cells.Item(2, 15).Value = 0.005520525738044089
$ws.Cells.Item(2, 16).Value = 0.005624540846623205
$ws.Cells.Item(2, 17).Value = 25.79038268582122
$ws.Cells.Item(2, 18).Value = 232.113444172391
$ws.Cells.Item(2, 19).Value = 0.0007419989698536164
$ws.Cells.Item(2, 20).Value = 0.0007836003069049704
$ws.Cells.Item(3, 7).Value = 8.872289666666667
$ws.Cells.Item(3, 8).Value = 26.616869
$ws.Cells.Item(3, 9).Value = 0.1344073019604298
$ws.Cells.Item(3, 10).Value = 0.1393180933827548
$ws.Cells.Item(3, 15).Value = 0.3528665483720876
$ws.Cells.Item(3, 16).Value = 0.3595150912979765
$ws.Cells.Item(3, 17).Value = 1648.495768586942
$ws.Cells.Item(3, 18).Value = 14836.46191728248
$ws.Cells.Item(3, 19).Value = 0.04742784071878179
$ws.Cells.Item(3, 20).Value = 0.0500869570619611
$ws.Cells.Item(4, 7).Value = 8.872289666666667
$ws.Cells.Item(4, 8).Value = 26.616869
$ws.Cells.Item(4, 9).Value = 0.1344073019604298
$ws.Cells.Item(4, 10).Value = 0.1393180933827548
$ws.Cells.Item(4, 13).Value = 137.0717086666666
$ws.Cells.Item(4, 14).Value = 411.2151259999999
$ws.Cells.Item(4, 15).Value = 0.2603191943704447
$ws.Cells.Item(4, 16).Value = 0.2652240042658267
$ws.Cells.Item(4, 17).Value = 1216.13990439561
$ws.Cells.Item(4, 18).Value = 10945.25913956049
$ws.Cells.Item(4, 19).Value = 0.03498880056384417
$ws.Cells.Item(4, 20).Value = 0.03695050259365459
$ws.Cells.Item(5, 7).Value = 8.872289666666667
$ws.Cells.Item(5, 8).Value = 26.616869
$ws.Cells.Item(5, 9).Value = 0.1344073019604298
$ws.Cells.Item(5, 10).Value = 0.1393180933827548
$ws.Cells.Item(5, 13).Value = 29.2127365
$ws.Cells.Item(5, 14).Value = 58.425473
$ws.Cells.Item(5, 15).Value = 0.05547925319534149
$ws.Cells.Item(5, 16).Value = 0.03768304451958546
$ws.Cells.Item(5, 17).Value = 259.1838601840062
$ws.Cells.Item(5, 18).Value = 1555.103161104037
$ws.Cells.Item(5, 19).Value = 0.007456816736765403
$ws.Cells.Item(5, 20).Value = 0.005249929915326112
$ws.Cells.Item(6, 7).Value = 8.872289666666667
$ws.Cells.Item(6, 8).Value = 26.616869
$ws.Cells.Item(6, 9).Value = 0.1344073019604298
$ws.Cells.Item(6, 10).Value = 0.1393180933827548
$ws.Cells.Item(6, 13).Value = 171.5584106666666
$ws.Cells.Item(6, 14).Value = 514.6752319999999
$ws.Cells.Item(6, 15).Value = 0.3258144783240821
$ws.Cells.Item(6, 16).Value = 0.331953319069988
$ws.Cells.Item(6, 17).Value = 1522.115914187623
$ws.Cells.Item(6, 18).Value = 13699.04322768861
$ws.Cells.Item(6, 19).Value = 0.04379184497118482
$ws.Cells.Item(6, 20).Value = 0.04624710350490799
$ws.Cells.Item(7, 9).Value = 0.756455981800989
$ws.Cells.Item(7, 10).Value = 0.784094342906462
$ws.Cells.Item(7, 13).Value = 2.906846333333333
$ws.Cells.Item(7, 14).Value = 8.720538999999999
$ws.Cells.Item(7, 15).Value = 0.005520525738044089
$ws.Cells.Item(7, 16).Value = 0.005624540846623205
$ws.Cells.Item(7, 17).Value = 145.1505161629519
$ws.Cells.Item(7, 18).Value = 1306.354645466567
$ws.Cells.Item(7, 19).Value = 0.004176034717229771
$ws.Cells.Item(7, 20).Value = 0.004410170659283578
$ws.Cells.Item(8, 9).Value = 0.756455981800989
$ws.Cells.Item(8, 10).Value = 0.784094342906462
$ws.Cells.Item(8, 15).Value = 0.3528665483720876
$ws.Cells.Item(8, 16).Value = 0.3595150912979765
$ws.Cells.Item(8, 19).Value = 0.2669280112935338
$ws.Cells.Item(8, 20).Value = 0.2818937492762436
$ws.Cells.Item(9, 9).Value = 0.756455981800989
$ws.Cells.Item(9, 10).Value = 0.784094342906462
$ws.Cells.Item(9, 13).Value = 137.0717086666666
$ws.Cells.Item(9, 14).Value = 411.2151259999999
$ws.Cells.Item(9, 15).Value = 0.2603191943704447
$ws.Cells.Item(9, 16).Value = 0.2652240042658267
$ws.Cells.Item(9, 17).Value = 6844.541122161519
$ws.Cells.Item(9, 18).Value = 61600.87009945367
$ws.Cells.Item(9, 19).Value = 0.1969200117591372
$ws.Cells.Item(9, 20).Value = 0.2079606413478341
$ws.Cells.Item(10, 9).Value = 0.756455981800989
$ws.Cells.Item(10, 10).Value = 0.784094342906462
$ws.Cells.Item(10, 13).Value = 29.2127365
$ws.Cells.Item(10, 14).Value = 58.425473
$ws.Cells.Item(10, 15).Value = 0.05547925319534149
$ws.Cells.Item(10, 16).Value = 0.03768304451958546
$ws.Cells.Item(10, 17).Value = 1458.709300482678
$ws.Cells.Item(10, 18).Value = 8752.255802896068
$ws.Cells.Item(10, 19).Value = 0.04196761294546771
$ws.Cells.Item(10, 20).Value = 0.02954706203129931
$ws.Cells.Item(11, 9).Value = 0.756455981800989
$ws.Cells.Item(11, 10).Value = 0.784094342906462
$ws.Cells.Item(11, 13).Value = 171.5584106666666
$ws.Cells.Item(11, 14).Value = 514.6752319999999
$ws.Cells.Item(11, 15).Value = 0.3258144783240821
$ws.Cells.Item(11, 16).Value = 0.331953319069988
$ws.Cells.Item(11, 17).Value = 8566.600709094588
$ws.Cells.Item(11, 18).Value = 77099.40638185128
$ws.Cells.Item(11, 19).Value = 0.2464643110856206
$ws.Cells.Item(11, 20).Value = 0.2602827195918014
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.2238023333333334
$ws.Cells.Item(12, 8).Value = 0.6714070000000001
$ws.Cells.Item(12, 9).Value = 0.003390406414343712
$ws.Cells.Item(12, 10).Value = 0.003514280478437762
$ws.Cells.Item(12, 13).Value = 2.906846333333333
$ws.Cells.Item(12, 14).Value = 8.720538999999999
$ws.Cells.Item(12, 15).Value = 0.005520525738044089
$ws.Cells.Item(12, 16).Value = 0.005624540846623205
$ws.Cells.Item(12, 17).Value = 0.6505589920414444
$ws.Cells.Item(12, 18).Value = 5.855030928373
$ws.Cells.Item(12, 19).Value = 0.00001871682587281423
$ws.Cells.Item(12, 20).Value = 0.00001976621409746373
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.2238023333333334
$ws.Cells.Item(13, 8).Value = 0.6714070000000001
$ws.Cells.Item(13, 9).Value = 0.003390406414343712
$ws.Cells.Item(13, 10).Value = 0.003514280478437762
$ws.Cells.Item(13, 15).Value = 0.3528665483720876
$ws.Cells.Item(13, 16).Value = 0.3595150912979765
$ws.Cells.Item(13, 17).Value = 41.58308772153678
$ws.Cells.Item(13, 18).Value = 374.247789493831
$ws.Cells.Item(13, 19).Value = 0.001196361009008052
$ws.Cells.Item(13, 20).Value = 0.001263436867052249
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.2238023333333334
$ws.Cells.Item(14, 8).Value = 0.6714070000000001
$ws.Cells.Item(14, 9).Value = 0.003390406414343712
$ws.Cells.Item(14, 10).Value = 0.003514280478437762
$ws.Cells.Item(14, 13).Value = 137.0717086666666
$ws.Cells.Item(14, 14).Value = 411.2151259999999
$ws.Cells.Item(14, 15).Value = 0.2603191943704447
$ws.Cells.Item(14, 16).Value = 0.2652240042658267
$ws.Cells.Item(14, 17).Value = 30.67696823358688
$ws.Cells.Item(14, 18).Value = 276.092714102282
$ws.Cells.Item(14, 19).Value = 0.0008825878663703431
$ws.Cells.Item(14, 20).Value = 0.0009320715406044885
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.2238023333333334
$ws.Cells.Item(15, 8).Value = 0.6714070000000001
$ws.Cells.Item(15, 9).Value = 0.003390406414343712
$ws.Cells.Item(15, 10).Value = 0.003514280478437762
$ws.Cells.Item(15, 13).Value = 29.2127365
$ws.Cells.Item(15, 14).Value = 58.425473
$ws.Cells.Item(15, 15).Value = 0.05547925319534149
$ws.Cells.Item(15, 16).Value = 0.03768304451958546
$ws.Cells.Item(15, 17).Value = 6.537878591751833
$ws.Cells.Item(15, 18).Value = 39.22727155051101
$ws.Cells.Item(15, 19).Value = 0.0001880972158964847
$ws.Cells.Item(15, 20).Value = 0.0001324287877232803
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.2238023333333334
$ws.Cells.Item(16, 8).Value = 0.6714070000000001
$ws.Cells.Item(16, 9).Value = 0.003390406414343712
$ws.Cells.Item(16, 10).Value = 0.003514280478437762
$ws.Cells.Item(16, 13).Value = 171.5584106666666
$ws.Cells.Item(16, 14).Value = 514.6752319999999
$ws.Cells.Item(16, 15).Value = 0.3258144783240821
$ws.Cells.Item(16, 16).Value = 0.331953319069988
$ws.Cells.Item(16, 17).Value = 38.39517261015822
$ws.Cells.Item(16, 18).Value = 345.556553491424
$ws.Cells.Item(16, 19).Value = 0.001104643497196018
$ws.Cells.Item(16, 20).Value = 0.001166577068960281
$ws.Cells.Item(17, 7).Value = 6.980364
$ws.Cells.Item(17, 8).Value = 13.960728
$ws.Cells.Item(17, 9).Value = 0.1057463098242374
$ws.Cells.Item(17, 10).Value = 0.07307328323234559
$ws.Cells.Item(17, 13).Value = 2.906846333333333
$ws.Cells.Item(17, 14).Value = 8.720538999999999
$ws.Cells.Item(17, 15).Value = 0.005520525738044089
$ws.Cells.Item(17, 16).Value = 0.005624540846623205
$ws.Cells.Item(17, 17).Value = 20.290845498732
$ws.Cells.Item(17, 18).Value = 121.745072992392
$ws.Cells.Item(17, 19).Value = 0.0005837752250878872
$ws.Cells.Item(17, 20).Value = 0.0004110036663371944
$ws.Cells.Item(18, 7).Value = 6.980364
$ws.Cells.Item(18, 8).Value = 13.960728
$ws.Cells.Item(18, 9).Value = 0.1057463098242374
$ws.Cells.Item(18, 10).Value = 0.07307328323234559
$ws.Cells.Item(18, 15).Value = 0.3528665483720876
$ws.Cells.Item(18, 16).Value = 0.3595150912979765
$ws.Cells.Item(18, 17).Value = 1296.970787645604
$ws.Cells.Item(18, 18).Value = 7781.824725873624
$ws.Cells.Item(18, 19).Value = 0.03731433535076404
$ws.Cells.Item(18, 20).Value = 0.02627094809271962
$ws.Cells.Item(19, 7).Value = 6.980364
$ws.Cells.Item(19, 8).Value = 13.960728
$ws.Cells.Item(19, 9).Value = 0.1057463098242374
$ws.Cells.Item(19, 10).Value = 0.07307328323234559
$ws.Cells.Item(19, 13).Value = 137.0717086666666
$ws.Cells.Item(19, 14).Value = 411.2151259999999
$ws.Cells.Item(19, 15).Value = 0.2603191943704447
$ws.Cells.Item(19, 16).Value = 0.2652240042658267
$ws.Cells.Item(19, 17).Value = 956.8104205952877
$ws.Cells.Item(19, 18).Value = 5740.862523571727
$ws.Cells.Item(19, 19).Value = 0.02752779418109292
$ws.Cells.Item(19, 20).Value = 0.01938078878373359
$ws.Cells.Item(20, 7).Value = 6.980364
$ws.Cells.Item(20, 8).Value = 13.960728
$ws.Cells.Item(20, 9).Value = 0.1057463098242374
$ws.Cells.Item(20, 10).Value = 0.07307328323234559
$ws.Cells.Item(20, 13).Value = 29.2127365
$ws.Cells.Item(20, 14).Value = 58.425473
$ws.Cells.Item(20, 15).Value = 0.05547925319534149
$ws.Cells.Item(20, 16).Value = 0.03768304451958546
$ws.Cells.Item(20, 17).Value = 203.915534206086
$ws.Cells.Item(20, 18).Value = 815.662136824344
$ws.Cells.Item(20, 19).Value = 0.005866726297211895
$ws.Cells.Item(20, 20).Value = 0.002753623785236756
$ws.Cells.Item(21, 7).Value = 6.980364
$ws.Cells.Item(21, 8).Value = 13.960728
$ws.Cells.Item(21, 9).Value = 0.1057463098242374
$ws.Cells.Item(21, 10).Value = 0.07307328323234559
$ws.Cells.Item(21, 13).Value = 171.5584106666666
$ws.Cells.Item(21, 14).Value = 514.6752319999999
$ws.Cells.Item(21, 15).Value = 0.3258144783240821
$ws.Cells.Item(21, 16).Value = 0.331953319069988
$ws.Cells.Item(21, 17).Value = 1197.5401537148155
$ws.Cells.Item(21, 18).Value = 7185.240922288895
$ws.Cells.Item(21, 19).Value = 0.03445367877008067
$ws.Cells.Item(21, 20).Value = 0.024256918904318418
